$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price strings so exact formatting is preserved
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = '25.812.62'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = '1.636.31'
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = '215.42'
$ws.Range("E5").Value = '  +0.04%  '
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = '0.0642'
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("D10").Value = '19.76'
$ws.Range("E10").Value = '  -0.85%  '
$ws.Range("D11").Value = '0.0780'
$ws.Range("E11").Value = '  +0.04%  '
$ws.Range("D12").Value = '1.670.50'
$ws.Range("E12").Value = '  +2.01%  '
$ws.Range("E13").Value = '  -0.65%  '
$ws.Range("D14").Value = '1.860.55'
$ws.Range("D15").Value = '0.556'
$ws.Range("E15").Value = '  -0.64%  '
$ws.Range("D16").Value = '0.0₃0776'
$ws.Range("E16").Value = '  +1.89%  '
$ws.Range("D17").Value = '63.25'
$ws.Range("E17").Value = '  +0.41%  '
$ws.Range("D18").Value = '25.813.02'
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("E19").Value = '  -0.11%  '
$ws.Range("E20").Value = '  +2.67%  '
$ws.Range("D21").Value = '194.51'
$ws.Range("E21").Value = '  -0.28%  '
$ws.Range("E22").Value = '  +0.92%  '
$ws.Range("E23").Value = '  +0.88%  '
$ws.Range("D24").Value = '1.01'
$ws.Range("E24").Value = '  +0.19%  '
$ws.Range("E25").Value = '  -0.83%  '
$ws.Range("D26").Value = '139.73'
$ws.Range("E26").Value = '  -0.47%  '
$ws.Range("E27").Value = '  -3.85%  '
$ws.Range("E28").Value = '  +0.44%  '
$ws.Range("D29").Value = '15.66'
$ws.Range("E29").Value = '  +1.53%  '
$ws.Range("E30").Value = '  +0.19%  '
$ws.Range("D31").Value = '0.0491'
$ws.Range("E31").Value = '  +0.25%  '
$ws.Range("D32").Value = '3.34'
$ws.Range("E32").Value = '  +1.45%  '
$ws.Range("E33").Value = '  +1.38%  '
$ws.Range("E34").Value = '  +1.70%  '
$ws.Range("E35").Value = '  +0.35%  '
$ws.Range("D36").Value = '0.899'
$ws.Range("E36").Value = '  -0.53%  '
$ws.Range("E37").Value = '  +0.24%  '
$ws.Range("D38").Value = '0.553'
$ws.Range("E38").Value = '  +0.18%  '
$ws.Range("D39").Value = '1.109.18'
$ws.Range("E39").Value = '  -1.62%  '
$ws.Range("E40").Value = '  +0.48%  '
$ws.Range("E41").Value = '  +0.45%  '
$ws.Range("D42").Value = '5.58'
$ws.Range("E42").Value = '  +0.70%  '
$ws.Range("D43").Value = '0.805'
$ws.Range("E43").Value = '  +0.68%  '
$ws.Range("D44").Value = '99.28'
$ws.Range("E44").Value = '  +1.30%  '
$ws.Range("E45").Value = '  -4.23%  '
$ws.Range("D46").Value = '55.32'
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("D47").Value = '2.47'
$ws.Range("E47").Value = '  +12.43%  '
$ws.Range("D48").Value = '7.71'
$ws.Range("E48").Value = '  -0.54%  '
$ws.Range("E49").Value = '  -2.00%  '
$ws.Range("E50").Value = '  +0.07%  '
$ws.Range("E51").Value = '  -0.01%  '
